# Apply updated crypto price/volume data to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "28.061.97"
$ws.Range("E2").Value = "  -1.65%  "

# Row 3
$ws.Range("D3").Value = "1.792.33"
$ws.Range("E3").Value = "  -1.82%  "

# Row 4
$ws.Range("D4").Value = "'1.004"
$ws.Range("E4").Value = "  +0.31%  "

# Row 5
$ws.Range("D5").Value = "'313.43"
$ws.Range("E5").Value = "  -0.71%  "

# Row 7
$ws.Range("D7").Value = "'0.5189"
$ws.Range("E7").Value = "  +1.44%  "

# Row 8
$ws.Range("D8").Value = "'0.3812"
$ws.Range("E8").Value = "  -3.68%  "

# Row 9
$ws.Range("D9").Value = "'0.07835"
$ws.Range("E9").Value = "  -4.67%  "

# Row 10
$ws.Range("D10").Value = "'41.36"
$ws.Range("E10").Value = "  -0.82%  "

# Row 11
$ws.Range("D11").Value = "'1.094"
$ws.Range("E11").Value = "  -1.78%  "

# Row 12
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").Value = "'6.263"
$ws.Range("E12").Value = "  -1.28%  "

# Row 13
$ws.Range("B13").Value = "BinanceUSD"
$ws.Range("C13").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D13").Value = "'1.004"
$ws.Range("E13").Value = "  +0.30%  "

# Row 14
$ws.Range("D14").Value = "'20.47"
$ws.Range("E14").Value = "  -3.37%  "

# Row 15
$ws.Range("D15").Value = "1.797.42"
$ws.Range("E15").Value = "  -1.59%  "

# Row 16
$ws.Range("D16").Value = "'7.269"
$ws.Range("E16").Value = "  -3.76%  "

# Row 17
$ws.Range("D17").Value = "'92.06"
$ws.Range("E17").Value = "  -0.90%  "

# Row 18
$ws.Range("D18").Value = "'0.00001080"
$ws.Range("E18").Value = "  -4.20%  "

# Row 19
$ws.Range("D19").Value = "'0.06539"
$ws.Range("E19").Value = "  -1.79%  "

# Row 20
$ws.Range("D20").Value = "'1.003"
$ws.Range("E20").Value = "  +0.31%  "

# Row 21
$ws.Range("D21").Value = "'17.27"
$ws.Range("E21").Value = "  -3.22%  "

# Row 22
$ws.Range("D22").Value = "'5.938"
$ws.Range("E22").Value = "  -2.58%  "

# Row 23
$ws.Range("D23").Value = "28.105.84"
$ws.Range("E23").Value = "  -1.62%  "

# Row 24
$ws.Range("D24").Value = "'11.10"
$ws.Range("E24").Value = "  -2.93%  "

# Row 25
$ws.Range("D25").Value = "'2.257"
$ws.Range("E25").Value = "  -0.20%  "

# Row 26
$ws.Range("D26").Value = "'160.78"
$ws.Range("E26").Value = "  +2.51%  "

# Row 27
$ws.Range("D27").Value = "'20.39"
$ws.Range("E27").Value = "  -4.67%  "

# Row 28
$ws.Range("D28").Value = "1.993.95"
$ws.Range("E28").Value = "  -2.16%  "

# Row 29
$ws.Range("D29").Value = "'2.322"
$ws.Range("E29").Value = "  -3.71%  "

# Row 30
$ws.Range("D30").Value = "'122.57"
$ws.Range("E30").Value = "  -3.63%  "

# Row 31
$ws.Range("D31").Value = "'0.1062"
$ws.Range("E31").Value = "  -2.60%  "

# Row 32
$ws.Range("D32").Value = "'1.048"
$ws.Range("E32").Value = "  -5.80%  "

# Row 33
$ws.Range("D33").Value = "'3.673"
$ws.Range("E33").Value = "  +0.47%  "

# Row 34
$ws.Range("D34").Value = "'5.525"
$ws.Range("E34").Value = "  -4.20%  "

# Row 35
$ws.Range("D35").Value = "'0.07249"
$ws.Range("E35").Value = "  +2.46%  "

# Row 36
$ws.Range("D36").Value = "'12.17"
$ws.Range("E36").Value = "  +7.79%  "

# Row 37
$ws.Range("D37").Value = "'0.02312"
$ws.Range("E37").Value = "  -1.78%  "

# Row 38
$ws.Range("D38").Value = "'8.729"
$ws.Range("E38").Value = "  -0.87%  "

# Row 39
$ws.Range("D39").Value = "'0.2130"
$ws.Range("E39").Value = "  -4.59%  "

# Row 40
$ws.Range("D40").Value = "'5.059"
$ws.Range("E40").Value = "  -4.14%  "

# Row 41
$ws.Range("D41").Value = "'0.6128"
$ws.Range("E41").Value = "  -3.19%  "

# Row 42
$ws.Range("E42").Value = "  -2.17%  "

# Row 43
$ws.Range("D43").Value = "'1.369"
$ws.Range("E43").Value = "  -2.12%  "

# Row 44
$ws.Range("D44").Value = "'13.23"
$ws.Range("E44").Value = "  -2.29%  "

# Row 45
$ws.Range("D45").Value = "'3.766"
$ws.Range("E45").Value = "  +0.84%  "

# Row 46
$ws.Range("D46").Value = "'0.5912"
$ws.Range("E46").Value = "  -0.57%  "

# Row 47
$ws.Range("D47").Value = "'127.71"
$ws.Range("E47").Value = "  +1.97%  "

# Row 48
$ws.Range("D48").Value = "'1.232"
$ws.Range("E48").Value = "  +3.22%  "

# Row 49
$ws.Range("D49").Value = "'1.911"
$ws.Range("E49").Value = "  -4.31%  "

# Row 50
$ws.Range("D50").Value = "'0.06729"
$ws.Range("E50").Value = "  -3.10%  "

# Row 51
$ws.Range("D51").Value = "'72.75"
$ws.Range("E51").Value = "  -1.90%  "
